$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-18
# from serial date 45172 (2023-09-03) to 45175 (2023-09-06)
$newDate = Get-Date -Year 2023 -Month 9 -Day 6 -Hour 0 -Minute 0 -Second 0
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate.Date
}
